$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the status/transaction/date fields for row 2 to reflect a new
# successful run (PASSED) with its transaction id and timestamp.
$ws.Range("G2").Value = "PASSED"
$ws.Range("H2").Value = "AAACT23179Z2PQC82"
$ws.Range("I2").Value = "28 jun. 2023, 14:43:14"

$wb.Save()
